# Insert a new data row before the current row 100, shifting existing
# rows 100-159 down to 101-160, then populate the new row 100 with the
# new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 100 (pushes old row 100 -> 101, etc.)
$ws.Rows.Item(100).Insert()

# Populate the new row 100 with the new record's values.
$ws.Cells.Item(100, 1).Value = 9
$ws.Cells.Item(100, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(100, 3).Value = 'Metropolitana'
$ws.Cells.Item(100, 4).Value = 45072
$ws.Cells.Item(100, 5).Value = 13
$ws.Cells.Item(100, 6).Value = 100112022
$ws.Cells.Item(100, 7).Value = 'Arveja Verde'
$ws.Cells.Item(100, 8).Value = 'Perfection'
$ws.Cells.Item(100, 9).Value = 'Primera'
$ws.Cells.Item(100, 10).Value = 43
$ws.Cells.Item(100, 11).Value = 27000
$ws.Cells.Item(100, 12).Value = 30000
$ws.Cells.Item(100, 13).Value = 28535
$ws.Cells.Item(100, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(100, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(100, 16).Value = 1141
$ws.Cells.Item(100, 17).Value = 25
$ws.Cells.Item(100, 18).Value = 'Hortaliza'
